$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.199468085106383
$ws.Range("C2").Value = 0.5079787234042553
$ws.Range("J2").Value = 0.02393617021276596
$ws.Range("P2").Value = 0.1595744680851064
$ws.Range("S2").Value = 0.1090425531914894
$ws.Range("B3").Value = 0.01015228426395939
$ws.Range("C3").Value = 0.02538071065989848
$ws.Range("J3").Value = 0.03553299492385787
$ws.Range("P3").Value = 0.7208121827411168
$ws.Range("S3").Value = 0.2081218274111675
$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.7297297297297297
$ws.Range("S4").Value = 0.2432432432432433
$ws.Range("B6").Value = 0.1062992125984252
$ws.Range("D6").Value = 0.01181102362204724
$ws.Range("F6").Value = 0.08267716535433071
$ws.Range("J6").Value = 0.1850393700787402
$ws.Range("O6").Value = 0.03149606299212598
$ws.Range("Q6").Value = 0.1259842519685039
$ws.Range("R6").Value = 0.04330708661417323
$ws.Range("S6").Value = 0.4133858267716535
$ws.Range("B7").Value = 0.1050228310502283
$ws.Range("D7").Value = 0.0136986301369863
$ws.Range("F7").Value = 0.0684931506849315
$ws.Range("J7").Value = 0.1141552511415525
$ws.Range("O7").Value = 0.0045662100456621
$ws.Range("Q7").Value = 0.0958904109589041
$ws.Range("R7").Value = 0.0776255707762557
$ws.Range("S7").Value = 0.5205479452054794
$ws.Range("B8").Value = 0.1297709923664122
$ws.Range("D8").Value = 0.01017811704834606
$ws.Range("E8").Value = 0.005089058524173028
$ws.Range("F8").Value = 0.06106870229007633
$ws.Range("J8").Value = 0.08651399491094147
$ws.Range("O8").Value = 0.02544529262086514
$ws.Range("Q8").Value = 0.2010178117048346
$ws.Range("R8").Value = 0.07124681933842239
$ws.Range("S8").Value = 0.4096692111959288
$ws.Range("B9").Value = 0.124223602484472
$ws.Range("D9").Value = 0.02484472049689441
$ws.Range("F9").Value = 0.08695652173913043
$ws.Range("J9").Value = 0.124223602484472
$ws.Range("O9").Value = 0.03105590062111801
$ws.Range("Q9").Value = 0.1863354037267081
$ws.Range("R9").Value = 0.06832298136645963
$ws.Range("S9").Value = 0.3540372670807453
$ws.Range("B10").Value = 0.1386217948717949
$ws.Range("D10").Value = 0.01923076923076923
$ws.Range("E10").Value = 0.0008012820512820513
$ws.Range("F10").Value = 0.08814102564102565
$ws.Range("J10").Value = 0.1233974358974359
$ws.Range("O10").Value = 0.01602564102564102
$ws.Range("Q10").Value = 0.1762820512820513
$ws.Range("R10").Value = 0.07051282051282051
$ws.Range("S10").Value = 0.3669871794871795
$ws.Range("G11").Value = 0.1589595375722543
$ws.Range("J11").Value = 0.08670520231213873
$ws.Range("K11").Value = 0.1936416184971098
$ws.Range("L11").Value = 0.5433526011560693
$ws.Range("S11").Value = 0.01734104046242774
$ws.Range("G12").Value = 0.7371134020618557
$ws.Range("J12").Value = 0.1907216494845361
$ws.Range("K12").Value = 0.02061855670103093
$ws.Range("L12").Value = 0.0154639175257732
$ws.Range("S12").Value = 0.03608247422680412
$ws.Range("G13").Value = 0.6595744680851063
$ws.Range("J13").Value = 0.2765957446808511
$ws.Range("S13").Value = 0.06382978723404255
$ws.Range("F15").Value = 0.02109704641350211
$ws.Range("H15").Value = 0.109704641350211
$ws.Range("I15").Value = 0.08016877637130802
$ws.Range("J15").Value = 0.3544303797468354
$ws.Range("K15").Value = 0.06751054852320675
$ws.Range("M15").Value = 0.02109704641350211
$ws.Range("O15").Value = 0.0759493670886076
$ws.Range("S15").Value = 0.270042194092827
$ws.Range("F16").Value = 0.008928571428571428
$ws.Range("H16").Value = 0.15625
$ws.Range("I16").Value = 0.0625
$ws.Range("J16").Value = 0.4330357142857143
$ws.Range("K16").Value = 0.1294642857142857
$ws.Range("M16").Value = 0.004464285714285714
$ws.Range("N16").Value = 0.004464285714285714
$ws.Range("O16").Value = 0.0625
$ws.Range("S16").Value = 0.1383928571428572
$ws.Range("F17").Value = 0.02638522427440633
$ws.Range("H17").Value = 0.1345646437994723
$ws.Range("I17").Value = 0.08970976253298153
$ws.Range("J17").Value = 0.4036939313984169
$ws.Range("K17").Value = 0.1240105540897098
$ws.Range("M17").Value = 0.0158311345646438
$ws.Range("O17").Value = 0.04749340369393139
$ws.Range("S17").Value = 0.158311345646438
$ws.Range("F18").Value = 0.006369426751592357
$ws.Range("H18").Value = 0.1210191082802548
$ws.Range("I18").Value = 0.09554140127388536
$ws.Range("J18").Value = 0.4585987261146497
$ws.Range("K18").Value = 0.1337579617834395
$ws.Range("M18").Value = 0.01273885350318471
$ws.Range("O18").Value = 0.06369426751592357
$ws.Range("S18").Value = 0.1082802547770701
$ws.Range("F19").Value = 0.01553254437869822
$ws.Range("H19").Value = 0.1982248520710059
$ws.Range("I19").Value = 0.05843195266272189
$ws.Range("J19").Value = 0.356508875739645
$ws.Range("K19").Value = 0.1161242603550296
$ws.Range("M19").Value = 0.02514792899408284
$ws.Range("O19").Value = 0.08431952662721894
$ws.Range("S19").Value = 0.1457100591715976
